$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.467.94"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "1.859.20"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'311.76"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.4768"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +3.00%  "
$ws.Range("D9").Value = "'0.07314"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "'0.9308"
$ws.Range("D11").Value = "'20.66"
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("D12").Value = "'0.07791"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.859.05"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "'5.452"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "'6.557"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "'90.09"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "'0.000008804"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "27.443.83"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "'14.62"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'10.69"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").Value = "'1.942"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "'154.76"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").Value = "'18.45"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").Value = "'2.005"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").Value = "'115.28"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "'4.939"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").Value = "'0.08877"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").Value = "'1.204"
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("D33").Value = "'0.7521"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "'4.575"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").Value = "'2.701"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'0.02041"
$ws.Range("E36").Value = "  +3.96%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "'0.5552"
$ws.Range("E38").Value = "  +5.19%  "
$ws.Range("D39").Value = "'0.05272"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").Value = "'7.006"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "'8.558"
$ws.Range("E42").Value = "  +3.17%  "
$ws.Range("D43").Value = "'0.1516"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.4865"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'10.63"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").Value = "'1.010"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").Value = "'103.46"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "'1.661"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("D49").Value = "'67.28"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("D50").Value = "'0.06091"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'0.9117"
$ws.Range("E51").Value = "  +2.19%  "
